# Apply the xfold-shortcuts.xlsx edit:
#  - shortcuts-default sheet: move "xfold-cmd" text from B4 to B2,
#    rewrite the A/B "快捷键/功能" pairs for rows 31-38, append 5 new
#    blank-but-formatted rows (46-50), and update the view/selection.
#  - shortcuts-leexioua sheet: update the view/selection only.

$wb = $excel.ActiveWorkbook

$wsDefault  = $wb.Worksheets.Item("shortcuts-default")
$wsLeexioua = $wb.Worksheets.Item("shortcuts-leexioua")

# ---------------------------------------------------------------------
# 1. shortcuts-default: B2 / B4 swap ("xfold-cmd" moves from the
#    "长按" (long-press) row to the "双键" (double-key) row).
# ---------------------------------------------------------------------
$wsDefault.Range("B2").Value = "xfold-cmd"
$wsDefault.Range("B4").ClearContents()

# ---------------------------------------------------------------------
# 2. shortcuts-default: rewrite rows 31-38, column A (快捷键) and
#    column B (功能/说明). Values only - C/D stay blank.
# ---------------------------------------------------------------------
$wsDefault.Range("A31").Value = "F1"
$wsDefault.Range("B31").Value = "xfold-cmd"

$wsDefault.Range("A32").Value = "``"
$wsDefault.Range("B32").Value = "xfold-cmd"

$wsDefault.Range("A33").Value = "Ctrl+Alt+R"
$wsDefault.Range("B33").Value = "重启xfold"

$wsDefault.Range("A34").Value = "Ctrl+Alt+\"
$wsDefault.Range("B34").Value = "SwitchOn开关切换"

$wsDefault.Range("A35").Value = "Ctrl+Q Q "
$wsDefault.Range("B35").Value = "关闭窗口"

$wsDefault.Range("A36").Value = "Win+-"
$wsDefault.Range("B36").Value = "窗口透明度-"

$wsDefault.Range("A37").Value = "Win+="
$wsDefault.Range("B37").Value = "窗口透明度+"

$wsDefault.Range("A38").Value = "LButton+Ctrl+C"
$wsDefault.Range("B38").Value = "复制文本（移除前置符号）"

# ---------------------------------------------------------------------
# 3. shortcuts-default: B32 / B34 swap their "shrink to fit" styling
#    (B32 loses it, B34 gains it) - copy the cell format, not just
#    the value.
# ---------------------------------------------------------------------
$wsDefault.Range("B32").Copy()
$wsDefault.Range("B34").PasteSpecial(-4122)

$wsDefault.Range("B33").Copy()
$wsDefault.Range("B32").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. shortcuts-default: append 5 new rows (46-50) that copy the
#    formatting of row 45 (content stays empty).
# ---------------------------------------------------------------------
for ($r = 46; $r -le 50; $r++) {
    $wsDefault.Cells.Item($r, 1).Value = "x"
    $wsDefault.Cells.Item($r, 1).ClearContents()
}

$wsDefault.Range("A45:D45").Copy()
$wsDefault.Range("A46:D50").PasteSpecial(-4122)

$wsDefault.Range("F45:M45").Copy()
$wsDefault.Range("F46:M50").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5. View / selection updates.
# ---------------------------------------------------------------------
$win1 = $wsDefault.Application.ActiveWindow
$wsDefault.Activate()
$win1.ScrollRow = 1
$win1.ScrollColumn = 1
[void]$wsDefault.Range("C18").Select()

$wsLeexioua.Activate()
$win2 = $wsLeexioua.Application.ActiveWindow
$win2.ScrollRow = 1
$win2.ScrollColumn = 1
[void]$wsLeexioua.Range("C28").Select()

$wsDefault.Activate()
